$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 11 (RE09): center the "NUM" column vertically, and vertically
# center the description column B as well (wraps already applied) ---
$ws.Range("A11").WrapText = $true
$ws.Range("A11").VerticalAlignment = -4108
$ws.Range("B11").VerticalAlignment = -4108

# --- Row 13 (RE11): new requirement about login/password access control,
# linked to the new "MANTER USUARIOS / AUTENTICAR USUARIO" use cases ---
$ws.Rows.Item(13).RowHeight = 30
$ws.Range("A13").WrapText = $true
$ws.Range("A13").VerticalAlignment = -4108
$ws.Range("B13").Value = "O acesso ao sistema deverá ser controlado através de login e senha"
$ws.Range("B13").VerticalAlignment = -4108
$ws.Range("C13").Value = "MANTER USUÁRIOS                AUTENTICAR USUÁRIO"
$ws.Range("C13").Font.Underline = $false

# --- Row 34 (UC08): new use case "AUTENTICAR USUARIO" with its actors ---
$ws.Range("B34").Value = "AUTENTICAR USUÁRIO"
$ws.Range("B34").WrapText = $true
$ws.Range("B34").VerticalAlignment = -4160
$ws.Range("C34").Value = "MASTER, ASSISTENTE, GERENTE, ANALISTA"

# --- Restore the selection that was active when the workbook was saved ---
$ws.Range("C35").Select()
